# Swap the two theme colour palettes that live in this deck.
#
# Before the edit:
#   ppt/theme/theme1.xml (the Slide Master's theme)  -> "Integral" palette
#   ppt/theme/theme2.xml (the Notes Master's theme)   -> "Office Theme" palette
#
# After the edit (per the authoritative OOXML diff) the two theme parts
# trade clrScheme contents (the fontScheme/fmtScheme blocks were already
# byte-identical in both themes, so a colour swap is the entire delta):
#   theme1.xml -> "Office Theme" palette
#   theme2.xml -> "Integral" palette
#
# The only colour-editing surface this host's PowerPoint object model
# exposes is Master.ColorScheme (the classic 12-slot scheme: dk1, lt1,
# dk2, lt2, accent1-6, hlink, folHlink, addressed as Colors(1..12) with
# an RGB value packed as 0x00BBGGRR). That maps onto
# $p.SlideMaster -> ppt/theme/theme1.xml, so we drive the Slide Master's
# theme to the target "Office Theme" palette through it.

$p = $ppt.ActivePresentation

function Set-RGB($colorScheme, [int]$index, [string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    $colorScheme.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

# Target palette = the "Office Theme" colours (dk1/lt1 stay black/white).
$officeTheme = @{
    1  = "000000"  # dk1
    2  = "FFFFFF"  # lt1
    3  = "44546A"  # dk2
    4  = "E7E6E6"  # lt2
    5  = "5B9BD5"  # accent1
    6  = "ED7D31"  # accent2
    7  = "A5A5A5"  # accent3
    8  = "FFC000"  # accent4
    9  = "4472C4"  # accent5
    10 = "70AD47"  # accent6
    11 = "0563C1"  # hlink
    12 = "954F72"  # folHlink
}

$masterColors = $p.SlideMaster.ColorScheme
foreach ($idx in $officeTheme.Keys) {
    Set-RGB $masterColors $idx $officeTheme[$idx]
}
